# issue #5: property land done
# Normalises the "土地" (land) sheet to the common schema used by the other
# sheets (English field headers + property_category/category/date/
# legislator_name/legislator_id/source_file/index columns), and cleans up
# stray whitespace / dashes that had crept into a handful of shared strings
# across every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "土地" (land) - rewrite headers to the English field names and add
# the trailing metadata columns (I:O), matching every other sheet.
# ---------------------------------------------------------------------
$land = $wb.Worksheets.Item("土地")

$land.Range("B1").Value = "name"
$land.Range("C1").Value = "area"
$land.Range("D1").Value = "share_portion"
$land.Range("E1").Value = "owner"
$land.Range("F1").Value = "register_date"
$land.Range("G1").Value = "register_reason"
$land.Range("H1").Value = "acquire_value"

$land.Range("I1").Value = "property_category"
$land.Range("J1").Value = "category"
$land.Range("K1").Value = "date"
$land.Range("L1").Value = "legislator_name"
$land.Range("M1").Value = "legislator_id"
$land.Range("N1").Value = "source_file"
$land.Range("O1").Value = "index"

# header style (bold + border) for the newly added header cells
$land.Range("B1").Copy()
$land.Range("I1:O1").PasteSpecial(-4122)

# clean up the stray inner spaces / dashes in the existing land strings
$land.Range("B2").Value = "臺北市士林區陽明段四小段06180000地號"
$land.Range("D2").Value = "10000分之121"
$land.Range("F2").Value = "75年01月21日"

$land.Range("B3").Value = "新北市中和區景平段06560000地號"
$land.Range("F3").Value = "81年06月18日"
$land.Range("G3").Value = "地籍圖重測"

# new metadata columns, row 2 (index 14)
$land.Range("I2").Value = "land"
$land.Range("J2").Value = "normal"
$land.Range("K2").NumberFormat = "@"
$land.Range("K2").Value = "2012-04-27"
$land.Range("L2").Value = "陳雪生"
$land.Range("M2").Value = 1751
$land.Range("N2").Value = "tmp5a001"
$land.Range("O2").Value = 14

# new metadata columns, row 3 (index 15)
$land.Range("I3").Value = "land"
$land.Range("J3").Value = "normal"
$land.Range("K3").NumberFormat = "@"
$land.Range("K3").Value = "2012-04-27"
$land.Range("L3").Value = "陳雪生"
$land.Range("M3").Value = 1751
$land.Range("N3").Value = "tmp5a001"
$land.Range("O3").Value = 15

# re-apply the plain data style (no special "text" number format) to the
# two date cells so they fall back in line with their row neighbours
$land.Range("B2").Copy()
$land.Range("K2").PasteSpecial(-4122)
$land.Range("B3").Copy()
$land.Range("K3").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Sheet "建物" (building) - strip stray inner spaces / dashes
# ---------------------------------------------------------------------
$building = $wb.Worksheets.Item("建物")
$building.Range("B2").Value = "臺北市士林區陽明段四小段40893000建號"
$building.Range("F2").Value = "75年01月21曰"
$building.Range("H2").Value = "(超過五年含陽台11.78平方公尺）"
$building.Range("B3").Value = "臺北市士林區陽明段四小段40937000建號"
$building.Range("D3").Value = "10000分之98"
$building.Range("F3").Value = "75年01月21日"
$building.Range("B4").Value = "新北市中和區景平段00310000建號"
$building.Range("F4").Value = "62年06月05日"
$building.Range("G4").Value = "第一次登記"

# ---------------------------------------------------------------------
# Sheet "汽車" (car)
# ---------------------------------------------------------------------
$car = $wb.Worksheets.Item("汽車")
$car.Range("E2").Value = "88年10月16日"

# ---------------------------------------------------------------------
# Sheet "存款" (deposits)
# ---------------------------------------------------------------------
$deposit = $wb.Worksheets.Item("存款")
$deposit.Range("B2").Value = "台北富邦商業銀行士林分行"
$deposit.Range("B7").Value = "中華郵政股份有限公司馬祖郵局"
$deposit.Range("B8").Value = "中華郵政股份有限公司馬袓郵局"
$deposit.Range("B14").Value = "台新國際商業銀行天母分行"
$deposit.Range("B15").Value = "台新國際商業銀行天母分行"
$deposit.Range("B16").Value = "台新國際商業銀行天母分行"
$deposit.Range("B20").Value = "台北富邦商業銀行承德分行"
$deposit.Range("B21").Value = "台北富邦商業銀行士林分行"
$deposit.Range("B22").Value = "台北富邦商業銀行士林分行"
$deposit.Range("B23").Value = "台北富邦商業銀行士林分行"

# ---------------------------------------------------------------------
# Sheet "基金受益憑證" (funds)
# ---------------------------------------------------------------------
$fund = $wb.Worksheets.Item("基金受益憑證")
$fund.Range("D2").Value = "曰盛國際商業銀行"
$fund.Range("D5").Value = "台北富邦商業銀行承德分行"
$fund.Range("H5").Value = "t65312"
$fund.Range("D6").Value = "台北富邦商業銀行承德分行"
$fund.Range("D7").Value = "台北富邦商業銀行承德分行"
$fund.Range("D8").Value = "台北富邦商業銀行承德分行"

# ---------------------------------------------------------------------
# Sheet "保險" (insurance)
# ---------------------------------------------------------------------
$insurance = $wb.Worksheets.Item("保險")
$insurance.Range("C1").Value = "保險名稱"
$insurance.Range("C7").Value = "大多利率變動型年金保險(甲型）"
$insurance.Range("C8").Value = "心得意利率變動型年金保險(甲型）"
$insurance.Range("E12").Value = "1000923投保"
$insurance.Range("C13").Value = "美利成增外幣養老保險（美元）"
$insurance.Range("E13").Value = "1000608投保"
$insurance.Range("E14").Value = "1010101投保"
$insurance.Range("E19").Value = "1000704投保"
